# Weekly update: insert two new daily price records for "Perejil" at
# Vega Modelo de Temuco, pushing the existing historical rows down by
# two positions (rows 132-191 -> 134-193) and populating the freshly
# opened rows 132 and 133 with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 132 (shifts old 132.. down by 2).
$ws.Rows.Item(132).Insert()
$ws.Rows.Item(132).Insert()

# --- New row 132 ---
$ws.Cells.Item(132, 1).Value = 10
$ws.Cells.Item(132, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(132, 3).Value = "La Araucanía"
$ws.Cells.Item(132, 4).Value = 44466
$ws.Cells.Item(132, 5).Value = 9
$ws.Cells.Item(132, 6).Value = 100112044
$ws.Cells.Item(132, 7).Value = "Perejil"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 90
$ws.Cells.Item(132, 11).Value = 3000
$ws.Cells.Item(132, 12).Value = 4000
$ws.Cells.Item(132, 13).Value = 3556
$ws.Cells.Item(132, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(132, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(132, 16).Value = 1185
$ws.Cells.Item(132, 17).Value = 3
$ws.Cells.Item(132, 18).Value = "Hortaliza"

# --- New row 133 ---
$ws.Cells.Item(133, 1).Value = 10
$ws.Cells.Item(133, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(133, 3).Value = "La Araucanía"
$ws.Cells.Item(133, 4).Value = 44466
$ws.Cells.Item(133, 5).Value = 9
$ws.Cells.Item(133, 6).Value = 100112044
$ws.Cells.Item(133, 7).Value = "Perejil"
$ws.Cells.Item(133, 8).Value = "Sin especificar"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 20
$ws.Cells.Item(133, 11).Value = 3300
$ws.Cells.Item(133, 12).Value = 3300
$ws.Cells.Item(133, 13).Value = 3300
$ws.Cells.Item(133, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(133, 15).Value = "Región Metropolitana"
$ws.Cells.Item(133, 16).Value = 1100
$ws.Cells.Item(133, 17).Value = 3
$ws.Cells.Item(133, 18).Value = "Hortaliza"
